# Commit: "adding delt kpi, extra only own manuf"
#
# The Assortment sheet's row 2 had two cells (B2, C2) both holding the
# shared string "General". These two cells are removed entirely (not
# just cleared to blank) so that the row collapses back to only A2.
# Because "General" then becomes an orphan entry in the shared string
# table, it disappears from xl/sharedStrings.xml once the workbook is
# re-serialized.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assortment")

# Remove the contents of B2:C2 so the underlying <c> elements are
# dropped from the sheet XML (matching the diff, which deletes the
# cells rather than leaving them blank).
$ws.Range("B2:C2").ClearContents()
